$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Add the new sheet positioned right after MNIST_Network, matching the target sheet order/sheetId.
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "CIFAR - 10"

# --- Cell values (row-major order so new shared strings are appended in the right sequence) ---
$ws2.Range("B2").Value = "Accuracy"
$ws2.Range("C2").Value = "Epoch"
$ws2.Range("D2").Value = "MaxPoolLayers"
$ws2.Range("E2").Value = "Conv2DLayers"
$ws2.Range("F2").Value = "Dense"
$ws2.Range("G2").Value = "Strides"
$ws2.Range("H2").Value = "Padding"
$ws2.Range("I2").Value = "Trainable Parameters"
$ws2.Range("J2").Value = "optimizer"
$ws2.Range("K2").Value = "activation"
$ws2.Range("A3").Value = "Architecture-1"
$ws2.Range("B3").Value = 67
$ws2.Range("C3").Value = 50
$ws2.Range("D3").Value = 2
$ws2.Range("E3").Value = 7
$ws2.Range("F3").Value = 2
$ws2.Range("G3").Value = 1
$ws2.Range("H3").Value = 1
$ws2.Range("I3").Value = 305162
$ws2.Range("J3").Value = "rmsprop"
$ws2.Range("K3").Value = "relu"
$ws2.Range("A4").Value = "Architecture-2"
$ws2.Range("B4").Value = 95.18
$ws2.Range("C4").Value = 50
$ws2.Range("D4").Value = 2
$ws2.Range("E4").Value = 8
$ws2.Range("F4").Value = 2
$ws2.Range("G4").Value = 1
$ws2.Range("H4").Value = 2
$ws2.Range("I4").Value = 314410
$ws2.Range("J4").Value = "adam"
$ws2.Range("K4").Value = "relu"
$ws2.Range("A5").Value = "Architecture-3"
$ws2.Range("B5").Value = 92.77
$ws2.Range("C5").Value = 50
$ws2.Range("D5").Value = 2
$ws2.Range("E5").Value = 8
$ws2.Range("F5").Value = 2
$ws2.Range("G5").Value = 1
$ws2.Range("H5").Value = 2
$ws2.Range("I5").Value = 314410
$ws2.Range("J5").Value = "adam"
$ws2.Range("K5").Value = "LeakyReLU"
$ws2.Range("A6").Value = "Architecture-4"
$ws2.Range("B6").Value = 90.24
$ws2.Range("C6").Value = 50
$ws2.Range("D6").Value = 2
$ws2.Range("E6").Value = 7
$ws2.Range("F6").Value = 2
$ws2.Range("G6").Value = 1
$ws2.Range("H6").Value = 4
$ws2.Range("I6").Value = 162442
$ws2.Range("J6").Value = "rmsprop learnign rate 0.0001"
$ws2.Range("K6").Value = "LeakyReLU"
$ws2.Range("A7").Value = "Architecture-5"
$ws2.Range("B7").Value = 92.29
$ws2.Range("C7").Value = 50
$ws2.Range("D7").Value = 2
$ws2.Range("E7").Value = 8
$ws2.Range("F7").Value = 2
$ws2.Range("G7").Value = 1
$ws2.Range("H7").Value = 2
$ws2.Range("I7").Value = 314410
$ws2.Range("J7").Value = "adam"
$ws2.Range("K7").Value = "LeakyReLU"

# --- Styles: reuse existing cellXfs where possible, create the two new ones by cloning + tweaking ---
# style index 2 (body, fontId=2 Arial/theme color) donor: sheet1!I1
foreach ($r in @("H2","I2","J2","K2","A3","B3","C3","D3","E3","F3","G3","H3","J3","K3","A4","B4","C4","D4","E4","F4","G4","H4","J4","K4","A5","B5","C5","D5","E5","F5","G5","H5","J5","K5","A6","B6","D6","E6","F6","G6","H6","J6","K6","A7","B7","C7","D7","E7","F7","G7","H7","J7","K7")) {
    $ws1.Range("I1").Copy()
    $ws2.Range($r).PasteSpecial(-4122)
}

# style index 3 (fontId=1) donor: sheet1!F2
foreach ($r in @("C6")) {
    $ws1.Range("F2").Copy()
    $ws2.Range($r).PasteSpecial(-4122)
}

# style index 4 (fontId=1, numFmt #,##0) donor: sheet1!J2
foreach ($r in @("I3")) {
    $ws1.Range("J2").Copy()
    $ws2.Range($r).PasteSpecial(-4122)
}

# style index 5 (fontId=2, centered) is new -> clone the fontId=2 look, then center it, then fan out
$ws1.Range("I1").Copy()
$ws2.Range("B2").PasteSpecial(-4122)
$ws2.Range("B2").HorizontalAlignment = -4108
foreach ($r in @("C2","D2","E2","F2","G2")) {
    $ws2.Range("B2").Copy()
    $ws2.Range($r).PasteSpecial(-4122)
}

# style index 6 (fontId=2, numFmt #,##0) is new -> clone the fontId=2 look, then apply number format, then fan out
$ws1.Range("I1").Copy()
$ws2.Range("I4").PasteSpecial(-4122)
$ws2.Range("I4").NumberFormat = "#,##0"
foreach ($r in @("I5","I6","I7")) {
    $ws2.Range("I4").Copy()
    $ws2.Range($r).PasteSpecial(-4122)
}

# --- Column I width (target stores 19.43 post-rounding; 18.6 lands on the nearest achievable snap) ---
$ws2.Columns.Item(9).ColumnWidth = 18.6

Write-Output "CIFAR - 10 sheet added"
